# Expanded types of violence covered by VESSA
#
# Applies the three text edits from the commit:
#   1. "... of domestic and/or sexual violence and for employees who have a
#      family or household member who is a ..." becomes
#      "... of domestic violence, sexual violence, gender-based violence,
#      stalking, or other violence, and for employees who have a family or
#      household member who is a ..."
#   2. ", and who claim protection under that law." becomes
#      " and who claim protection under that law." (drop the comma)
#   3. "... of domestic or sexual violence or because ... resulting from the
#      domestic or sexual violen..." becomes
#      "... of violence or because ... resulting from the violen..."

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    " of domestic and/or sexual violence and for employees who have a family or household member who is a ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " of domestic violence, sexual violence, gender-based violence, stalking, or other violence, and for employees who have a family or household member who is a ",
    2)

$d.Content.Find.Execute(
    ", and who claim protection under that law.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " and who claim protection under that law.",
    2)

$d.Content.Find.Execute(
    " of domestic or sexual violence or because they requested an adjustment to their work schedule to cope with the violence and its many effects. VESSA also specifically requires that an employer reasonably accommodate any known work-related limitations resulting from the domestic or sexual violen",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " of violence or because they requested an adjustment to their work schedule to cope with the violence and its many effects. VESSA also specifically requires that an employer reasonably accommodate any known work-related limitations resulting from the violen",
    2)
